$wb = $excel.ActiveWorkbook

# Update "展览" sheet (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5136
$ws1.Range("F3").Value = 156
$ws1.Range("F4").Value = 902

# Update "全部类型" sheet (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5136
$ws4.Range("F3").Value = 156
$ws4.Range("F4").Value = 902
